$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("제출답안")

$ws.Range("A4").Value = "a"
$ws.Range("B4").Value = "v"
$ws.Range("C4").Value = "c"

$ws.Range("G3").Select()
